# Carretta & Enriquez 2012b - add effort estimates sheet
# 1) Rename existing "Sheet1" -> "Estimates"
# 2) Add a new "Effort" sheet right after "Estimates"
# 3) Populate the Effort sheet with headers (bold) + one data row
# 4) Autofit columns, set selections to match final UI state

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet ---
$wsEstimates = $wb.Worksheets.Item(1)
$wsEstimates.Name = "Estimates"

# --- Add the new "Effort" sheet directly after "Estimates" ---
$wsEffort = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsEstimates)
$wsEffort.Name = "Effort"

# --- Headers (row 1), bold ---
$wsEffort.Range("A1").Value = "year"
$wsEffort.Range("B1").Value = "nvessels"
$wsEffort.Range("C1").Value = "mesh_in_avg"
$wsEffort.Range("D1").Value = "sets_tot_est"
$wsEffort.Range("E1").Value = "sets_obs"
$wsEffort.Range("F1").Value = "obs_perc"
$wsEffort.Range("A1:F1").Font.Bold = $true

# --- Data (row 2) ---
$wsEffort.Range("A2").Value = 2011
$wsEffort.Range("B2").Value = 50
$wsEffort.Range("C2").Value = 7.4
$wsEffort.Range("D2").Value = 2123
$wsEffort.Range("E2").Value = 171
$wsEffort.Range("F2").Value = 8

# --- Size columns to fit their content, like the source workbook's other sheet ---
$wsEffort.Columns("A:F").AutoFit() | Out-Null

# --- Restore/assign the on-screen selections for each sheet ---
$wsEstimates.Range("C35").Select() | Out-Null

$wsEffort.Select() | Out-Null
$wsEffort.Range("A1:F1").Select() | Out-Null

Write-Output "Added Effort sheet with effort estimates."
